# Auto-generated edit script for LOB1264.xlsx
# Reproduces the commit "Build site at 2022-09-26 16:07:08 UTC" transformation
# applied to docs/assets/disciplinas/LOB1264.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target content/layout for rows 1-24 (cols A-C) after the edit.
$rows = @(
  @{ r=1; ht=$null; cells=@{ B="Ementa atual:"; C="Ementa modificada (dados modificados em vermelho):" } },
  @{ r=2; ht=$null; cells=@{ B="LOB1264"; C="LOB1264" } },
  @{ r=3; ht=$null; cells=@{ A="Nome:"; B=" Toxicologia Ambiental"; C=" Toxicologia Ambiental" } },
  @{ r=4; ht=$null; cells=@{ A="Name:"; B="Environmental Toxicology"; C="Environmental Toxicology" } },
  @{ r=5; ht=$null; cells=@{ A="Créditos-aula:"; B="2"; C="2" } },
  @{ r=6; ht=$null; cells=@{ A="Créditos-trabalho"; B="0"; C="0" } },
  @{ r=7; ht=$null; cells=@{ A="Carga horária:"; B="30 h"; C="30 h" } },
  @{ r=8; ht=$null; cells=@{ A="Ativação:"; B="01/01/2022"; C="01/01/2022" } },
  @{ r=9; ht=$null; cells=@{ A="Semestre ideal:"; B="EA-4"; C="EA-4" } },
  @{ r=10; ht=60.0; cells=@{ A="Objetivos:"; B="8855158 - Morun Bernardino Neto"; C="8855158 - Morun Bernardino Neto" } },
  @{ r=11; ht=60.0; cells=@{ A="Objectives:"; B="To train higher education professionals with knowledge in environmental toxicology that will allow them to act in the assessment of environmental risk and public health, assessment of environmental impact and public health, as well as to act in the management of this risk."; C="To train higher education professionals with knowledge in environmental toxicology that will allow them to act in the assessment of environmental risk and public health, assessment of environmental impact and public health, as well as to act in the management of this risk." } },
  @{ r=12; ht=$null; cells=@{ A="Docentes responsáveis:" } },
  @{ r=13; ht=60.0; cells=@{ A="Programa resumido:"; B="Semestral"; C="Semestral" } },
  @{ r=14; ht=60.0; cells=@{ A="Short syllabus:"; B="General principles and molecular bases of the mechanisms of toxicity; Toxicokinetics; Toxic agents of environmental interest and their effects"; C="General principles and molecular bases of the mechanisms of toxicity; Toxicokinetics; Toxic agents of environmental interest and their effects" } },
  @{ r=15; ht=120.0; cells=@{ A="Programa:"; B="01/01/2022"; C="01/01/2022" } },
  @{ r=16; ht=120.0; cells=@{ A="Syllabus:"; B="General principles of toxicology: Principles of toxicology; Molecular basis of the mechanisms of toxicity. Toxicokinetics: Absorption; Bioavailability and biotransformation; Apparent distribution volume; Clearance; Half-life period; Elimination; Dose-dependent toxicokinetics; Accumulation during continuous or intermittent exposure. Toxic Agents: Toxic effects of pesticides; Toxic effects of solvents; Toxic effects of metals; Toxic effects of poisons."; C="General principles of toxicology: Principles of toxicology; Molecular basis of the mechanisms of toxicity. Toxicokinetics: Absorption; Bioavailability and biotransformation; Apparent distribution volume; Clearance; Half-life period; Elimination; Dose-dependent toxicokinetics; Accumulation during continuous or intermittent exposure. Toxic Agents: Toxic effects of pesticides; Toxic effects of solvents; Toxic effects of metals; Toxic effects of poisons." } },
  @{ r=17; ht=$null; cells=@{ A="Avaliação:" } },
  @{ r=18; ht=60.0; cells=@{ A="Método:"; B="8855158 - Morun Bernardino Neto"; C="8855158 - Morun Bernardino Neto" } },
  @{ r=19; ht=60.0; cells=@{ A="Critério:"; B="Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos."; C="Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos." } },
  @{ r=20; ht=60.0; cells=@{ A="Norma de recuperação:"; B="O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2"; C="O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2" } },
  @{ r=21; ht=120.0; cells=@{ A="Bibliografia:"; B="Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2"; C="Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2" } },
  @{ r=22; ht=$null; cells=@{ A="Requisitos:" } },
  @{ r=23; ht=30.0; cells=@{ B="LOQ4081 -  Química Orgânica  (Requisito fraco)`n"; C="LOQ4081 -  Química Orgânica  (Requisito fraco)`n" } },
  @{ r=24; ht=30.0; cells=@{ B="LOT2046 -  Microbiologia e Bioquimica Aplicadas  (Requisito fraco)`n"; C="LOT2046 -  Microbiologia e Bioquimica Aplicadas  (Requisito fraco)`n" } }
)

foreach ($row in $rows) {
  $r = $row.r
  foreach ($col in @("A","B","C")) {
    if ($row.cells.ContainsKey($col)) {
      $ws.Range("$col$r").Value = $row.cells[$col]
    } else {
      $ws.Range("$col$r").ClearContents()
    }
  }
  if ($null -ne $row.ht) {
    $ws.Rows.Item($r).RowHeight = $row.ht
  }
}

# Row 25 no longer exists after the edit (dimension shrinks from C25 to C24).
$ws.Rows.Item(25).Delete()

